# Update countries & provincias Spain
# Applies the 26-Abril-2020 data refresh to the "Pais" sheet:
#  - updates the "last updated" timestamp
#  - refreshes case counts for a handful of countries (Australia, Oman,
#    Guinea Ecuatorial)
#  - a few countries tied on total cases change rank order, so the rows
#    around them are rewritten with the correct country name + stats

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp banner -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 26 de Abril de 2020 a las 09:52"

# --- Straight numeric refreshes (country/rank unchanged) --------------
# Australia (row 46)
$ws.Range("B46").Value = 6711
$ws.Range("C46").Value = 16
$ws.Range("D46").Value = 5539
$ws.Range("E46").Value = 1089

# Oman (row 67)
$ws.Range("B67").Value = 1998
$ws.Range("C67").Value = 93
$ws.Range("D67").Value = 333
$ws.Range("E67").Value = 1655

# Guinea Ecuatorial (row 127)
$ws.Range("D127").Value = 8
$ws.Range("E127").Value = 249

# --- Granada / Belice swap ranks (both tied on 18 total cases) --------
$ws.Range("A183").Value = "Granada"
$ws.Range("B183").Value = 18
$ws.Range("C183").Value = 0
$ws.Range("D183").Value = 7
$ws.Range("E183").Value = 11
$ws.Range("F183").Value = 4
$ws.Range("G183").Value = 0
$ws.Range("H183").Value = 0

$ws.Range("A184").Value = "Belice"
$ws.Range("B184").Value = 18
$ws.Range("C184").Value = 0
$ws.Range("D184").Value = 5
$ws.Range("E184").Value = 11
$ws.Range("F184").Value = 1
$ws.Range("G184").Value = 0
$ws.Range("H184").Value = 2

# --- Burundi / Islas Turcas y Caicos swap ranks (both tied on 11, -----
# --- identical stats so only the labels trade places) -----------------
$ws.Range("A198").Value = "Islas Turcas y Caicos"
$ws.Range("A199").Value = "Burundi"

# --- Sudan del Sur gets a new case and jumps the tied-at-6 group ------
$ws.Range("A208").Value = "Sudan del Sur"
$ws.Range("B208").Value = 6
$ws.Range("C208").Value = 1
$ws.Range("D208").Value = 0
$ws.Range("E208").Value = 6
$ws.Range("F208").Value = 0
$ws.Range("G208").Value = 0
$ws.Range("H208").Value = 0

$ws.Range("A209").Value = "Islas Virgenes Britanicas"
$ws.Range("B209").Value = 6
$ws.Range("C209").Value = 0
$ws.Range("D209").Value = 3
$ws.Range("E209").Value = 2
$ws.Range("F209").Value = 0
$ws.Range("G209").Value = 0
$ws.Range("H209").Value = 1

$ws.Range("A210").Value = "Sahara Occidental"
$ws.Range("B210").Value = 6
$ws.Range("C210").Value = 0
$ws.Range("D210").Value = 5
$ws.Range("E210").Value = 1
$ws.Range("F210").Value = 0
$ws.Range("G210").Value = 0
$ws.Range("H210").Value = 0

$ws.Range("A211").Value = "San Bartolome"
$ws.Range("B211").Value = 6
$ws.Range("C211").Value = 0
$ws.Range("D211").Value = 6
$ws.Range("E211").Value = 0
$ws.Range("F211").Value = 0
$ws.Range("G211").Value = 0
$ws.Range("H211").Value = 0

$ws.Range("A212").Value = "Bonaire, San Eustaquio y Saba"
$ws.Range("B212").Value = 5
$ws.Range("C212").Value = 0
$ws.Range("D212").Value = 0
$ws.Range("E212").Value = 5
$ws.Range("F212").Value = 0
$ws.Range("G212").Value = 0
$ws.Range("H212").Value = 0
